$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = ""
$ws.Range("A4").Value = ""
$ws.Range("A5").Value = ""
